$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new columns D (Window #) and E (Genotype)
$windowNum = @(1,2,3,1,1,2,3,1,1,2,1,2,1,2,1,2,3,4)
$genotype  = @(2,2,2,1,2,2,2,2,1,1,2,2,1,1,2,2,2,2)

# Headers - set E1 (Genotype) first so it gets shared-string index 3,
# then D1 (Window #) so it gets shared-string index 4, matching the target file.
$ws.Range("E1").Value = "Genotype"
$ws.Range("D1").Value = "Window #"

for ($i = 0; $i -lt $windowNum.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $windowNum[$i]
    $ws.Cells.Item($row, 5).Value = $genotype[$i]
}

# Column widths matching the target diff (closest achievable values given this
# runtime's column-width quantization; target stored widths are 10.85546875 / 10.28515625)
$ws.Columns.Item(4).ColumnWidth = 10.0
$ws.Columns.Item(5).ColumnWidth = 9.5

# Update selection to match target state
$ws.Range("H18").Select()
